$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.381.13"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.068.97"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.41"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.53"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.373.21"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.39"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.068.87"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.302.03"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.55"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.64"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.74"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("E28").Value = "  -6.46%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.09"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0957"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.72"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.480.76"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.15"
$ws.Range("E46").Value = "  -7.18%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.08"
$ws.Range("E49").Value = "  -5.51%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.259.89"
$ws.Range("E51").Value = "  +0.19%  "
